$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "Data Science 대학원 교수 채용 원칙"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/data-science-professor-interview/#utm_source=rss&utm_medium=rss&utm_campaign=data-science-professor-interview"

$ws.Range("D39").Value = "Probability concepts explained: probability distributions"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Probability-concepts-explained-probability-distributions"

$ws.Range("D41").Value = "Data Mesh – Part 2 사용시의 유의점"
$ws.Range("E41").Value = "http://cloudinsight.net/data/data-mesh-part-2-%ec%82%ac%ec%9a%a9%ec%8b%9c%ec%9d%98-%ec%9c%a0%ec%9d%98%ec%a0%90/"

$ws.Range("D46").Value = "[한국식품연구원] 2021년 03월, 생물정보학(Bioinformatics 채용), 식품연구데이터 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/386"

$ws.Range("D51").Value = "토트넘과 상대하는 GNK 디나모 자그레브란?"
$ws.Range("E51").Value = "https://bskyvision.com/1128"
